# Insert a new data row for "Vega Modelo de Temuco - Zanahoria" right after
# the existing row 255 (i.e. at row 256), pushing the previous rows 256-315
# down to 257-316, and populate the newly inserted row with its own values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 256..315 down to 257..316 by inserting a blank row at 256.
$ws.Rows.Item(256).Insert()

# Fill in the new row 256 with the new observation's data.
$ws.Cells.Item(256, 1).Value = 10
$ws.Cells.Item(256, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(256, 3).Value = "La Araucanía"
$ws.Cells.Item(256, 4).Value = 44736
$ws.Cells.Item(256, 5).Value = 9
$ws.Cells.Item(256, 6).Value = 100114013
$ws.Cells.Item(256, 7).Value = "Zanahoria"
$ws.Cells.Item(256, 8).Value = "Sin especificar"
$ws.Cells.Item(256, 9).Value = "Primera"
$ws.Cells.Item(256, 10).Value = 140
$ws.Cells.Item(256, 11).Value = 6000
$ws.Cells.Item(256, 12).Value = 6000
$ws.Cells.Item(256, 13).Value = 6000
$ws.Cells.Item(256, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(256, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(256, 16).Value = 240
$ws.Cells.Item(256, 17).Value = 25
$ws.Cells.Item(256, 18).Value = "Hortaliza"
